$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends at row 30, whose phone number is stored as the
# zero-padded text "09876543". The edit duplicates that record into a new
# row 31 (keeping the original zero-padded text phone number, blank
# birthday, and 0 points) and turns row 30's phone number into a plain
# numeric value (9876543) with points left at 0 -- i.e. "points
# 09876543 -> 0.00" becomes its own, newly numeric-keyed row.

# Push everything below row 30 down by one (there's nothing below it, so
# this simply opens up a new, blank row 31).
$ws.Rows.Item(31).Insert()

# Re-create row 30's original record in the new row 31. A leading
# apostrophe forces text storage so the leading zero on the phone number
# is preserved, and "'" alone yields an empty (but still text-typed)
# birthday cell, matching the untouched B30 cell it was copied from.
$ws.Cells.Item(31, 1).Value = "'09876543"
$ws.Cells.Item(31, 1).Style = "Normal"
$ws.Cells.Item(31, 2).Value = "'"
$ws.Cells.Item(31, 2).Style = "Normal"
$ws.Cells.Item(31, 3).Value = 0

# Row 30 keeps its blank birthday and 0 points, but the phone number
# becomes a genuine number (9876543) instead of zero-padded text.
$ws.Cells.Item(30, 1).Value = 9876543
